# Update "Forecast Comparison" sheet with corrected forecast output:
#  - insert a new "Week_Start_Date" column (B) between Week and ASIN
#  - rewrite the short week labels (W01 -> W1, etc.)
#  - fill in the Week_Start_Date column with the week's start date (as text)
#  - update the MyForecast values with corrected figures
#  - store the holiday flag as a proper boolean
# Also refresh the dependent numbers on the "Summary" sheet.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("Forecast Comparison")
$sum = $wb.Worksheets.Item("Summary")

# 1) Make room for the new "Week_Start_Date" column right after "Week".
$ws.Columns("B:B").Insert()

# 2) Header row.
$ws.Range("B1").Value = "Week_Start_Date"

# 3) Row data: Week, Week_Start_Date, MyForecast (corrected), is_holiday_week.
#    ASIN (C), Amazon Mean/P70/P80/P90 (E:H) and Product Title (I) are unchanged.
$weeks = @(
    @{ Row = 2;  Week = "W1";  Start = "2025-01-05"; Forecast = 170 },
    @{ Row = 3;  Week = "W2";  Start = "2025-01-12"; Forecast = 166 },
    @{ Row = 4;  Week = "W3";  Start = "2025-01-19"; Forecast = 162 },
    @{ Row = 5;  Week = "W4";  Start = "2025-01-26"; Forecast = 154 },
    @{ Row = 6;  Week = "W5";  Start = "2025-02-02"; Forecast = 151 },
    @{ Row = 7;  Week = "W6";  Start = "2025-02-09"; Forecast = 140 },
    @{ Row = 8;  Week = "W7";  Start = "2025-02-16"; Forecast = 134 },
    @{ Row = 9;  Week = "W8";  Start = "2025-02-23"; Forecast = 145 },
    @{ Row = 10; Week = "W9";  Start = "2025-03-02"; Forecast = 163 },
    @{ Row = 11; Week = "W10"; Start = "2025-03-09"; Forecast = 130 },
    @{ Row = 12; Week = "W11"; Start = "2025-03-16"; Forecast = 156 },
    @{ Row = 13; Week = "W12"; Start = "2025-03-23"; Forecast = 134 },
    @{ Row = 14; Week = "W13"; Start = "2025-03-30"; Forecast = 123 },
    @{ Row = 15; Week = "W14"; Start = "2025-04-06"; Forecast = 132 },
    @{ Row = 16; Week = "W15"; Start = "2025-04-13"; Forecast = 119 },
    @{ Row = 17; Week = "W16"; Start = "2025-04-20"; Forecast = 122 }
)

# The Week_Start_Date column must hold plain text (not be auto-converted to a
# date serial number), so force a text number format before assigning values.
$ws.Range("B2:B17").NumberFormat = "@"

foreach ($w in $weeks) {
    $r = $w.Row
    $ws.Cells.Item($r, 1).Value = $w.Week
    $ws.Cells.Item($r, 2).Value = $w.Start
    $ws.Cells.Item($r, 4).Value = $w.Forecast
    $ws.Cells.Item($r, 10).Value = $false
}

# 4) Summary sheet totals recomputed from the corrected forecast.
#    These values are stored as text (like the rest of the column), so keep
#    the cells formatted as text to avoid them being coerced into numbers.
$sum.Range("B9:B12").NumberFormat = "@"
$sum.Range("B9").Value  = "2300"
$sum.Range("B10").Value = "1223"
$sum.Range("B11").Value = "652"
$sum.Range("B12").Value = "170"
